# Apply updated optimization output values across several result sheets.
$wb = $excel.ActiveWorkbook

# --- optimized_production_rates (sheet10) ---
$ws = $wb.Worksheets.Item("optimized_production_rates")
$ws.Range("B2").Value = 0.30694938228558577
$ws.Range("B3").Value = 1.0743086107694393
$ws.Range("B4").Value = 0.48771524593590354

# --- optimized_threshold_b (sheet11) ---
$ws = $wb.Worksheets.Item("optimized_threshold_b")
$ws.Range("B2").Value = [double]"3.0587440664167575E-3"
$ws.Range("B3").Value = 1.3839158090808557
$ws.Range("B4").Value = 0.66309294236400151

# --- network_optimized_weights (sheet12) ---
$ws = $wb.Worksheets.Item("network_optimized_weights")
$ws.Range("C2").Value = 1.9939437371105837
$ws.Range("D3").Value = -1.9372125735418451
$ws.Range("B4").Value = 0.50183442381666499
$ws.Range("C4").Value = -0.51982030316205574

# --- optimization_diagnostics (sheet13) ---
$ws = $wb.Worksheets.Item("optimization_diagnostics")
$ws.Range("B2").Value = [double]"3.4717146051340112E-10"
$ws.Range("B3").Value = 1.2091802766687352
$ws.Range("B5").Value = 7104
$ws.Range("B8").Value = [double]"1.2993402247154986E-11"
$ws.Range("B9").Value = [double]"1.5689021843658784E-10"
$ws.Range("B10").Value = [double]"8.7163076085646052E-10"

# --- wt_log2_optimized_expression (sheet8) ---
$ws = $wb.Worksheets.Item("wt_log2_optimized_expression")
$ws.Range("C2").Value = 0.41646677733897119
$ws.Range("D2").Value = 0.46176807853492619
$ws.Range("E2").Value = 0.40160771727658495
$ws.Range("F2").Value = 0.32223209573017675
$ws.Range("G2").Value = 0.25476347511591979
$ws.Range("H2").Value = 0.20671453136009876
$ws.Range("I2").Value = 0.17607123791630941
$ws.Range("J2").Value = 0.15805653259498154
$ws.Range("K2").Value = 0.14817593846626842
$ws.Range("L2").Value = 0.14311546455371871
$ws.Range("M2").Value = 0.14072153460539172
$ws.Range("N2").Value = 0.1397090158556753

$ws.Range("C3").Value = -0.78490648797862961
$ws.Range("D3").Value = -1.4925198021243902
$ws.Range("E3").Value = -2.10723478868388
$ws.Range("F3").Value = -2.5916768687798313
$ws.Range("G3").Value = -2.9271336249489166
$ws.Range("H3").Value = -3.128805991991245
$ws.Range("I3").Value = -3.2343311476537009
$ws.Range("J3").Value = -3.2821213212851097
$ws.Range("K3").Value = -3.2999132449441522
$ws.Range("L3").Value = -3.3040597681564954
$ws.Range("M3").Value = -3.3029651882826796
$ws.Range("N3").Value = -3.3004761706499863

$ws.Range("C4").Value = 0.18266604131799169
$ws.Range("D4").Value = 0.370356878813632
$ws.Range("E4").Value = 0.48202682583066903
$ws.Range("F4").Value = 0.53599840082459982
$ws.Range("G4").Value = 0.5571379909751295
$ws.Range("H4").Value = 0.56226344292546881
$ws.Range("I4").Value = 0.56081237642383175
$ws.Range("J4").Value = 0.5574833466572664
$ws.Range("K4").Value = 0.55426743881401563
$ws.Range("L4").Value = 0.55179567203149005
$ws.Range("M4").Value = 0.55010840495467073
$ws.Range("N4").Value = 0.54904549961099991
